# Update "想去人数" (people interested) counts on both the "展览" sheet
# and the "全部类型" sheet (they mirror the same data). Column F holds the
# counts; rows are 1-indexed with header row 1.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 1987
    7  = 1654
    9  = 663
    19 = 3844
    24 = 711
    25 = 486
    26 = 354
    27 = 32
    28 = 1644
    30 = 158
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
